$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6 data - fill in this order so the shared-string table is built
# in the same order as the target (制版费, PCB制版费预付款, 嘉利创, 备注, then remark text)
$ws.Range("B6").Value = 20180407

$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "制版费"

$ws.Range("D6").Value = "PCB制版费预付款"
$ws.Range("D6").Characters(4,6).Font.Name = "宋体"

$ws.Range("C2").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "嘉利创"

# New header in H1 "备注" (Remarks) - copy format (style s="1") from an existing styled cell
$ws.Range("C2").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "备注"

$ws.Range("G6").Value = 202
$ws.Range("H6").Value = "充值200，1%手续费"
$ws.Range("H6").Characters(1,2).Font.Name = "宋体"
$ws.Range("H6").Characters(3,3).Font.Name = "Tahoma"
$ws.Range("H6").Characters(6,1).Font.Name = "宋体"
$ws.Range("H6").Characters(7,2).Font.Name = "Tahoma"
$ws.Range("H6").Characters(9,3).Font.Name = "宋体"

# New column width for column H (target stored width 24.625; engine snaps to
# pixel grid with MDW=7, so 23.91 lands on the closest reachable value)
$ws.Columns.Item(8).ColumnWidth = 23.91

# Update selection
$ws.Range("H11").Select()
